$d = $word.ActiveDocument

# Remove the stray leading "anan" text that was accidentally left in the
# first (bold, size-28) heading run, restoring it to just "1. Check your system".
$d.Content.Find.Execute("anan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
